$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '96.964.65'
$ws.Range("E2").Value = '  +2.93%  '

# Row 3
$ws.Range("D3").Value = '3.328.21'
$ws.Range("E3").Value = '  +7.14%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$helper.Value = '244.05'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +2.11%  '

# Row 6
$helper.Value = '623.44'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +1.25%  '

# Row 7
$helper.Value = '1.12'
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.60%  '

# Row 8
$helper.Value = '0.387'
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -0.77%  '

# Row 9
$ws.Range("E9").Value = '  -0.01%  '

# Row 10
$ws.Range("D10").Value = '3.329.75'
$ws.Range("E10").Value = '  +6.90%  '

# Row 11
$helper.Value = '0.787'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -4.54%  '

# Row 12
$helper.Value = '0.200'
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +0.64%  '

# Row 13
$ws.Range("D13").Value = '96.822.39'
$ws.Range("E13").Value = '  +3.09%  '

# Row 14
$helper.Value = '0.0000246'
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +0.63%  '

# Row 15
$helper.Value = '35.25'
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +1.81%  '

# Row 16
$ws.Range("D16").Value = '3.963.38'
$ws.Range("E16").Value = '  +7.24%  '

# Row 17
$helper.Value = '5.49'
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +1.87%  '

# Row 18
$ws.Range("D18").Value = '3.356.49'
$ws.Range("E18").Value = '  +6.27%  '

# Row 19
$helper.Value = '3.57'
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -3.92%  '

# Row 20
$helper.Value = '15.14'
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +1.56%  '

# Row 21
$helper.Value = '487.73'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  +8.64%  '

# Row 22
$helper.Value = '0.0000209'
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +3.85%  '

# Row 23
$helper.Value = '5.84'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -1.54%  '

# Row 24
$helper.Value = '9.26'
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +3.02%  '

# Row 25
$helper.Value = '5.67'
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +0.22%  '

# Row 26
$helper.Value = '88.55'
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +2.43%  '

# Row 27
$helper.Value = '12.09'
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +0.09%  '

# Row 28
$ws.Range("D28").Value = '3.537.41'
$ws.Range("E28").Value = '  +7.38%  '

# Row 29
$helper.Value = '0.999'
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +0.01%  '

# Row 30
$helper.Value = '0.180'
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -0.05%  '

# Row 31
$helper.Value = '0.241'
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -7.27%  '

# Row 32
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$helper.Value = '1.00'
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +0.04%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$helper.Value = '0.121'
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -3.65%  '

# Row 34
$helper.Value = '9.28'
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -0.11%  '

# Row 35
$helper.Value = '27.51'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  +5.20%  '

# Row 36
$helper.Value = '7.40'
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -6.58%  '

# Row 37
$helper.Value = '0.151'
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -6.66%  '

# Row 38
$helper.Value = '1.93'
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +0.93%  '

# Row 39
$helper.Value = '494.05'
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +3.16%  '

# Row 40
$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$helper.Value = '24.61'
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +2.78%  '

# Row 41
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$helper.Value = '0.449'
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -0.86%  '

# Row 42
$helper.Value = '1.27'
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -1.10%  '

# Row 43
$helper.Value = '0.803'
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +15.79%  '

# Row 44
$helper.Value = '3.26'
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  -1.17%  '

# Row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$helper.Value = '1.00'
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +0.06%  '

# Row 46
$ws.Range("B46").Value = 'MantraDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$helper.Value = '3.44'
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -8.02%  '

# Row 47
$helper.Value = '159.97'
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -0.26%  '

# Row 48
$helper.Value = '1.93'
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +3.54%  '

# Row 49
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$helper.Value = '4.54'
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  +2.40%  '

# Row 50
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$helper.Value = '45.16'
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +2.98%  '

# Row 51
$helper.Value = '1.35'
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +2.27%  '

$helper.Clear()
$excel.CutCopyMode = 0
